# Trade #105 closed at 2026-02-17 09:18:33 - unknown UNKNOWN +0.000%
#
# Updates the "Summary", "Strategy Status", "All Trades" and "MarketMaking"
# sheets to reflect the newly-closed trade #105 (row 106 in the trade logs).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Summary sheet - roll up totals
# ---------------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1200.12   # Current Capital
$summary.Range("B4").Value = 0.13      # Total P&L $
$summary.Range("B6").Value = 105       # Total Trades
$summary.Range("B7").Value = 45        # Winning Trades
$summary.Range("B9").Value = 42.86     # Win Rate %

# ---------------------------------------------------------------------------
# Strategy Status sheet - MarketMaking row (row 4)
# ---------------------------------------------------------------------------
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C4").Value = 100.12     # Capital
$status.Range("D4").Value = 105        # Trades
$status.Range("E4").Value = 0.13       # P&L $
$status.Range("F4").Value = 0.12       # P&L %
$status.Range("G4").Value = 42.86      # Win Rate %

# ---------------------------------------------------------------------------
# Append the new closed trade (row 106) to both "All Trades" and
# "MarketMaking" logs - the two sheets mirror the same trade history.
# ---------------------------------------------------------------------------
$sheetNames = @("All Trades", "MarketMaking")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)

    $ws.Cells.Item(106, 1).Value = 105                 # Trade #
    # Column B holds a literal "yyyy-mm-dd" date string, not a real date -
    # force Text so Excel does not auto-convert it to a date serial number.
    $ws.Cells.Item(106, 2).NumberFormat = "@"
    $ws.Cells.Item(106, 2).Value = "2026-02-17"        # Date
    $ws.Cells.Item(106, 3).Value = "09:18:27"          # Time
    $ws.Cells.Item(106, 4).Value = "MarketMaking"      # Strategy
    $ws.Cells.Item(106, 5).Value = "UP"                # Side
    $ws.Cells.Item(106, 6).Value = 0.16                # Entry Price
    $ws.Cells.Item(106, 7).Value = 0.18                # Exit Price
    $ws.Cells.Item(106, 8).Value = "CLOSED"            # Status
    $ws.Cells.Item(106, 9).Value = 12.5                # P&L %
    $ws.Cells.Item(106, 10).Value = 0.02               # P&L $
    $ws.Cells.Item(106, 11).Value = 100.12             # Capital After
    $ws.Cells.Item(106, 12).Value = 0                  # Entry Slippage (bps)
    $ws.Cells.Item(106, 13).Value = 0                  # Exit Slippage (bps)
    $ws.Cells.Item(106, 14).Value = 0.6                # Confidence
    $ws.Cells.Item(106, 15).Value = "Normal spread capture: 19600 bps"  # Entry Reason
    $ws.Cells.Item(106, 16).Value = "early_exit"       # Exit Reason
    $ws.Cells.Item(106, 17).Value = 0.13               # Duration (min)

    # Drop the Text-format override picked up above so the new row keeps the
    # same (default) styling as every other row in the sheet.
    $ws.Range("A106:Q106").ClearFormats()
}

Write-Host "Trade #105 appended; Summary and Strategy Status totals updated."
